$d = $word.ActiveDocument

$replacements = @(
    @{old = "170÷9=18, 8"; new = "770÷4=192, 2"},
    @{old = "981÷9=109, 0"; new = "591÷6=98, 3"},
    @{old = "753÷2=376, 1"; new = "524÷4=131, 0"},
    @{old = "948÷5=189, 3"; new = "864÷2=432, 0"},
    @{old = "820÷2=410, 0"; new = "834÷7=119, 1"},
    @{old = "374÷5=74, 4"; new = "305÷2=152, 1"},
    @{old = "798÷7=114, 0"; new = "371÷9=41, 2"},
    @{old = "627÷4=156, 3"; new = "975÷4=243, 3"},
    @{old = "979÷8=122, 3"; new = "735÷2=367, 1"},
    @{old = "322÷8=40, 2"; new = "864÷4=216, 0"},
    @{old = "755÷3=251, 2"; new = "674÷7=96, 2"},
    @{old = "328÷6=54, 4"; new = "480÷8=60, 0"},
    @{old = "253÷7=36, 1"; new = "830÷3=276, 2"},
    @{old = "835÷3=278, 1"; new = "723÷2=361, 1"},
    @{old = "854÷7=122, 0"; new = "165÷9=18, 3"},
    @{old = "604÷4=151, 0"; new = "225÷3=75, 0"},
    @{old = "917÷8=114, 5"; new = "693÷9=77, 0"},
    @{old = "710÷7=101, 3"; new = "144÷8=18, 0"},
    @{old = "382÷5=76, 2"; new = "209÷7=29, 6"},
    @{old = "790÷6=131, 4"; new = "232÷4=58, 0"},
    @{old = "252÷6=42, 0"; new = "599÷6=99, 5"},
    @{old = "431÷7=61, 4"; new = "759÷2=379, 1"},
    @{old = "908÷8=113, 4"; new = "842÷9=93, 5"},
    @{old = "492÷3=164, 0"; new = "911÷7=130, 1"},
    @{old = "422÷7=60, 2"; new = "206÷6=34, 2"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
